$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Pipeline(steps=[('scaler', None),`n                ('selector',`n                 <__main__.NamedFeatureSelector object at 0x7f181d355fa0>),`n                ('model',`n                 RandomForestClassifier(max_depth=2, max_features='log2',`n                                        min_samples_leaf=10,`n                                        min_samples_split=5, n_estimators=200,`n                                        random_state=42))])"
$ws.Range("B2").Value = 0.724044289044289
$ws.Range("C2").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7f181d319eb0>, 'scaler': None, 'model__n_estimators': 200, 'model__min_samples_split': 5, 'model__min_samples_leaf': 10, 'model__max_features': 'log2', 'model__max_depth': 2, 'model__class_weight': None}"
$ws.Range("D2").Value = 0.7908274173368717
$ws.Range("E2").Value = 0.6355930125430125
$ws.Range("F2").Value = 0.8333333333333334
$ws.Range("G2").Value = 0.7536348240209954
$ws.Range("H2").Value = 0.5975710317460318
$ws.Range("I2").Value = 0.75
$ws.Range("J2").Value = 0.8547234042553191
$ws.Range("K2").Value = 0.7133333333333332
$ws.Range("N2").Value = "[1 1 1 1 1 1 0 1 1 1 0 1 1 1 0 0 1 1 1 1 1 1 1 1]"
$ws.Range("A3").Value = "Pipeline(steps=[('scaler', MinMaxScaler()),`n                ('selector',`n                 <__main__.NamedFeatureSelector object at 0x7f181d3d18b0>),`n                ('model',`n                 RandomForestClassifier(max_depth=2, max_features='log2',`n                                        min_samples_leaf=11,`n                                        min_samples_split=7, n_estimators=200,`n                                        random_state=42))])"
$ws.Range("B3").Value = 0.7294971694971694
$ws.Range("C3").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7f181d20ea00>, 'scaler': MinMaxScaler(), 'model__n_estimators': 200, 'model__min_samples_split': 7, 'model__min_samples_leaf': 11, 'model__max_features': 'log2', 'model__max_depth': 2, 'model__class_weight': None}"
$ws.Range("D3").Value = 0.7912883190734448
$ws.Range("E3").Value = 0.5967024087024086
$ws.Range("F3").Value = 0.7894736842105263
$ws.Range("G3").Value = 0.7752558903792831
$ws.Range("H3").Value = 0.6825821428571427
$ws.Range("I3").Value = 0.6818181818181818
$ws.Range("J3").Value = 0.8264255319148937
$ws.Range("K3").Value = 0.5823333333333334
$ws.Range("L3").Value = 0.9375
$ws.Range("N3").Value = "[1 1 1 1 1 1 1 1 1 1 1 0 1 0 1 1 1 1 1 1 1 1 1 1]"
$ws.Range("A4").Value = "Pipeline(steps=[('scaler', None),`n                ('selector',`n                 <__main__.NamedFeatureSelector object at 0x7f181d3d1730>),`n                ('model',`n                 RandomForestClassifier(max_depth=1, max_features='log2',`n                                        min_samples_leaf=4,`n                                        min_samples_split=10,`n                                        random_state=42))])"
$ws.Range("B4").Value = 0.7045221445221446
$ws.Range("C4").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7f181d20eeb0>, 'scaler': None, 'model__n_estimators': 100, 'model__min_samples_split': 10, 'model__min_samples_leaf': 4, 'model__max_features': 'log2', 'model__max_depth': 1, 'model__class_weight': None}"
$ws.Range("D4").Value = 0.7963998176902769
$ws.Range("E4").Value = 0.6142599567099566
$ws.Range("F4").Value = 0.7647058823529413
$ws.Range("G4").Value = 0.7940489092838345
$ws.Range("H4").Value = 0.6364825396825396
$ws.Range("I4").Value = 0.8666666666666667
$ws.Range("J4").Value = 0.8065777777777778
$ws.Range("K4").Value = 0.6216
$ws.Range("L4").Value = 0.6842105263157895
$ws.Range("N4").Value = "[0 1 1 0 0 1 0 1 1 0 1 0 0 1 1 1 0 1 1 0 1 1 1 1]"
$ws.Range("A5").Value = "Pipeline(steps=[('scaler', RobustScaler()),`n                ('selector',`n                 <__main__.NamedFeatureSelector object at 0x7f181cf3d580>),`n                ('model',`n                 RandomForestClassifier(max_depth=1, min_samples_leaf=5,`n                                        min_samples_split=9, n_estimators=200,`n                                        random_state=42))])"
$ws.Range("B5").Value = 0.7690659340659339
$ws.Range("C5").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7f181d110b50>, 'scaler': RobustScaler(), 'model__n_estimators': 200, 'model__min_samples_split': 9, 'model__min_samples_leaf': 5, 'model__max_features': 'sqrt', 'model__max_depth': 1, 'model__class_weight': None}"
$ws.Range("D5").Value = 0.8030860187131779
$ws.Range("E5").Value = 0.6365851925851925
$ws.Range("F5").Value = 0.7368421052631579
$ws.Range("G5").Value = 0.7686470532440832
$ws.Range("H5").Value = 0.6159690476190476
$ws.Range("I5").Value = 0.5833333333333334
$ws.Range("J5").Value = 0.8629591836734692
$ws.Range("K5").Value = 0.7041666666666666
$ws.Range("N5").Value = "[1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1]"
$ws.Range("A6").Value = "Pipeline(steps=[('scaler', RobustScaler()),`n                ('selector',`n                 <__main__.NamedFeatureSelector object at 0x7f181d387280>),`n                ('model',`n                 RandomForestClassifier(max_depth=1, min_samples_leaf=11,`n                                        min_samples_split=9, n_estimators=50,`n                                        random_state=42))])"
$ws.Range("B6").Value = 0.772142857142857
$ws.Range("C6").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7f181d1b2940>, 'scaler': RobustScaler(), 'model__n_estimators': 50, 'model__min_samples_split': 9, 'model__min_samples_leaf': 11, 'model__max_features': 'sqrt', 'model__max_depth': 1, 'model__class_weight': None}"
$ws.Range("D6").Value = 0.8133185001371643
$ws.Range("E6").Value = 0.6690890054390054
$ws.Range("G6").Value = 0.7783157389766635
$ws.Range("H6").Value = 0.6239980158730157
$ws.Range("J6").Value = 0.8740192307692305
$ws.Range("K6").Value = 0.7518333333333334

Write-Output "Applied changes"
